$wb = $excel.ActiveWorkbook

# --- Sheet "MID_LFT_#1" : append row 87 ---
$ws = $wb.Worksheets.Item("MID_LFT_#1")
$ws.Range("A87").NumberFormat = $ws.Range("A86").NumberFormat
$ws.Range("A87").Value = 45873.46453703703
$ws.Range("B87").Value = "0x01,0x90"
$ws.Range("C87").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Range("D87").Value = "0x01,0x24"
$ws.Range("E87").Value = "0x07"
$ws.Range("F87").Value = 400
$ws.Range("G87").Value = 568631262647113000000000.0
$ws.Range("H87").Value = 292
$ws.Range("I87").Value = 7

# --- Sheet "MID_LFT_#2" : append row 87 ---
$ws = $wb.Worksheets.Item("MID_LFT_#2")
$ws.Range("A87").NumberFormat = $ws.Range("A86").NumberFormat
$ws.Range("A87").Value = 45873.46453703703
$ws.Range("B87").Value = "0x01,0x7c"
$ws.Range("C87").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Range("D87").Value = "0x01,0x2C"
$ws.Range("E87").Value = "0x19"
$ws.Range("F87").Value = 380
$ws.Range("G87").Value = 568432987514711000000000.0
$ws.Range("H87").Value = 300
$ws.Range("I87").Value = 25

# --- Sheet "MID_PLT_#1" : append row 87 ---
$ws = $wb.Worksheets.Item("MID_PLT_#1")
$ws.Range("A87").NumberFormat = $ws.Range("A86").NumberFormat
$ws.Range("A87").Value = 45873.46453703703
$ws.Range("B87").Value = "0x00,0x6e"
$ws.Range("C87").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Range("D87").Value = "0x00,0x5F"
$ws.Range("E87").Value = "0x15"
$ws.Range("F87").Value = 110
$ws.Range("G87").Value = 568631262647113000000000.0
$ws.Range("H87").Value = 95
$ws.Range("I87").Value = 15

# --- Sheet "MID_PLT_#2" : append row 87 ---
$ws = $wb.Worksheets.Item("MID_PLT_#2")
$ws.Range("A87").NumberFormat = $ws.Range("A86").NumberFormat
$ws.Range("A87").Value = 45873.46453703703
$ws.Range("B87").Value = "0x00,0x82"
$ws.Range("C87").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Range("D87").Value = "0x00,0x76"
$ws.Range("E87").Value = "0x9"
$ws.Range("F87").Value = 130
$ws.Range("G87").Value = 568631262647113000000000.0
$ws.Range("H87").Value = 118
$ws.Range("I87").Value = 9
